# Adds a new "2020-05-28" forecast column (AI) and a new "2020-05-28"
# observation row (47) to both the "cases" and "deaths" sheets, and fills
# in the newly observed value for 2020-05-14's forecast (B33).

$wb = $excel.ActiveWorkbook

$sheetNames = @("cases", "deaths")

# AI-column (column 35) values per sheet, keyed by row number.
# Row 1 is the header (reuses the existing "2020-05-14" shared string,
# following the sheet's existing pattern of header-date = row-date - 1 day).
# Rows 2-33 stay blank; rows 34-47 carry the new forecast figures.
$aiValues = @{
    "cases"  = @{
        1  = "2020-05-14"
        34 = 56048
        35 = 57410
        36 = 58969
        37 = 60450
        38 = 61984
        39 = 63093
        40 = 65018
        41 = 66323
        42 = 67574
        43 = 68724
        44 = 69833
        45 = 70778
        46 = 71910
        47 = 72998
    }
    "deaths" = @{
        1  = "2020-05-14"
        34 = 4427
        35 = 4507
        36 = 4595
        37 = 4669
        38 = 4748
        39 = 4825
        40 = 4898
        41 = 4975
        42 = 5053
        43 = 5123
        44 = 5191
        45 = 5259
        46 = 5326
        47 = 5394
    }
}

# New "observed" value landing in B33 (the 2020-05-14 row) for each sheet.
$b33Values = @{
    "cases"  = 54286
    "deaths" = 4315
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Track every cell we write a literal date-like string into, so we can
    # strip the auto-applied quote-prefix style back off at the end and
    # keep the cell on the sheet's default (unstyled) formatting.
    $textCells = @()

    # --- Column AI (35): header + the new forecast figures ---------------
    $rows = $aiValues[$sheetName].Keys | Sort-Object
    foreach ($row in $rows) {
        $val = $aiValues[$sheetName][$row]
        $cell = $ws.Cells.Item($row, 35)
        if ($row -eq 1) {
            # Force literal text so "2020-05-14" isn't reinterpreted as a
            # date serial number - matches the existing header cells.
            $cell.Value = "'" + $val
            $textCells += $cell
        } else {
            $cell.Value = $val
        }
    }

    # --- Row 47: new "2020-05-28" date label ------------------------------
    $a47 = $ws.Cells.Item(47, 1)
    $a47.Value = "'2020-05-28"
    $textCells += $a47

    # --- B33: newly observed value -----------------------------------------
    $ws.Range("B33").Value = $b33Values[$sheetName]

    # Strip the quote-prefix formatting picked up from the apostrophe-typed
    # text above, restoring the default (unstyled) cell format.
    foreach ($tc in $textCells) {
        $tc.ClearFormats()
    }
}
